$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new ErrorFolder row (row 5). Set A5/B5 text first so new shared
# strings are appended in the same order as the target workbook
# (ErrorFolder, then the Errors path, then the two new description strings).
$ws.Range("A5").Value = "ErrorFolder"
$ws.Range("B5").Value = "\\compass-usa\cgcorp\AccountingServices\Secure\Systems and Projects\Systems\SAP\Projects\UiPath\FB Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Errors"

# Add the new "Template where headers are pulled from." description for the
# existing TemplateFile row (row 4).
$ws.Range("C4").Value = "Template where headers are pulled from."

# Add the description for the new ErrorFolder row.
$ws.Range("C5").Value = "Where files are moved if there is an error reading or processing them."

# Turn B5 into a hyperlink pointing at the Errors folder (mirrors the other
# UNC-path hyperlinks already on the sheet).
$ws.Hyperlinks.Add($ws.Range("B5"), "file:///\\compass-usa\cgcorp\AccountingServices\Secure\Systems%20and%20Projects\Systems\SAP\Projects\UiPath\FB%20Dev\PriceAudits\Detail.5.9.2019.SC\Loaded\Errors")

# B5 should use the same Hyperlink cell style as the other Asset links (B2:B4).
$ws.Range("B5").Style = $ws.Range("B4").Style

# Move the view/selection down to the newly added row.
[void]$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollColumn = 3
